$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1912.375
$ws.Range("I18").Value = 1685.5714
$ws.Range("K18").Value = 1685.5714
$ws.Range("M18").Value = -1401.5714
$ws.Range("H40").Value = 3445.2727
$ws.Range("I40").Value = 3132
$ws.Range("K40").Value = 3132
$ws.Range("M40").Value = -2957
$ws.Range("H113").Value = 71432060
$ws.Range("I113").Value = 200003580
$ws.Range("J113").Value = 3444.4443
$ws.Range("K113").Value = 200003580
$ws.Range("L113").Value = 3444.4443
$ws.Range("M113").Value = -200000326
$ws.Range("N113").Value = -9952.444299999999
$ws.Range("H116").Value = 34385188
$ws.Range("J116").Value = 55562356
$ws.Range("L116").Value = 55562356
$ws.Range("N116").Value = -55569240
$ws.Range("H132").Value = 3255.5212
$ws.Range("I132").Value = 3180.8333
$ws.Range("J132").Value = 3411.3914
$ws.Range("K132").Value = 9542.499899999999
$ws.Range("L132").Value = 10234.1742
$ws.Range("M132").Value = -7012.499899999999
$ws.Range("N132").Value = -15294.1742
$ws.Range("H137").Value = 1949.8646
$ws.Range("I137").Value = 1304.1111
$ws.Range("J137").Value = 2337.3167
$ws.Range("K137").Value = 3912.3333
$ws.Range("L137").Value = 7011.9501
$ws.Range("M137").Value = -1362.3333
$ws.Range("N137").Value = -12111.9501
$ws.Range("H141").Value = 2922.4
$ws.Range("I141").Value = 2922.4
$ws.Range("K141").Value = 8767.200000000001
$ws.Range("M141").Value = -3587.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 180183.55
$ws.Range("I32").Value = 190158.22
$ws.Range("K32").Value = 190158.22
$ws.Range("M32").Value = -189871.22
$ws.Range("H74").Value = 4316.143
$ws.Range("I74").Value = 7000
$ws.Range("J74").Value = 2303.25
$ws.Range("K74").Value = 7000
$ws.Range("L74").Value = 2303.25
$ws.Range("M74").Value = -6126
$ws.Range("N74").Value = -4051.25
$ws.Range("H77").Value = 4316.143
$ws.Range("I77").Value = 7000
$ws.Range("J77").Value = 2303.25
$ws.Range("K77").Value = 35000
$ws.Range("L77").Value = 11516.25
$ws.Range("M77").Value = -30632
$ws.Range("N77").Value = -20252.25
$ws.Range("H110").Value = 76924240
$ws.Range("I110").Value = 76924240
$ws.Range("K110").Value = 76924240
$ws.Range("M110").Value = -76922195
$ws.Range("H122").Value = 3550.3333
$ws.Range("I122").Value = 1730
$ws.Range("J122").Value = 6644.9
$ws.Range("K122").Value = 5190
$ws.Range("L122").Value = 19934.7
$ws.Range("M122").Value = -2740
$ws.Range("N122").Value = -24834.7
$ws.Range("H132").Value = 2859616
$ws.Range("I132").Value = 2502624.5
$ws.Range("K132").Value = 7507873.5
$ws.Range("M132").Value = -7505343.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 381.66666
$ws.Range("J11").Value = 1699.6666
$ws.Range("L11").Value = 1699.6666
$ws.Range("N11").Value = -1979.6666
$ws.Range("H17").Value = 30
$ws.Range("I17").Value = 30
$ws.Range("J17").Value = 30
$ws.Range("K17").Value = 30
$ws.Range("L17").Value = 30
$ws.Range("M17").Value = 142
$ws.Range("N17").Value = -374
$ws.Range("H20").Value = 7438.8687
$ws.Range("I20").Value = 9838.23
$ws.Range("J20").Value = 2240.25
$ws.Range("K20").Value = 9838.23
$ws.Range("L20").Value = 2240.25
$ws.Range("M20").Value = -9591.23
$ws.Range("N20").Value = -2734.25
$ws.Range("H94").Value = 3242.1936
$ws.Range("I94").Value = 3174.5
$ws.Range("J94").Value = 3474.2856
$ws.Range("K94").Value = 3174.5
$ws.Range("L94").Value = 3474.2856
$ws.Range("M94").Value = -2723.5
$ws.Range("N94").Value = -4376.2856
$ws.Range("H105").Value = 2143.2917
$ws.Range("I105").Value = 1939.875
$ws.Range("K105").Value = 1939.875
$ws.Range("M105").Value = -192.875
$ws.Range("H107").Value = 16748601
$ws.Range("I107").Value = 106063.52
$ws.Range("J107").Value = 55581188
$ws.Range("K107").Value = 106063.52
$ws.Range("L107").Value = 55581188
$ws.Range("M107").Value = -104143.52
$ws.Range("N107").Value = -55585028

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 406.3125
$ws.Range("I7").Value = 366.73334
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 366.73334
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = -253.73334
$ws.Range("N7").Value = -1226
$ws.Range("H28").Value = 36166.332
$ws.Range("J28").Value = 36166.332
$ws.Range("L28").Value = 36166.332
$ws.Range("N28").Value = -36656.332
$ws.Range("H43").Value = 29999.5
$ws.Range("J43").Value = 29999.5
$ws.Range("L43").Value = 29999.5
$ws.Range("N43").Value = -30367.5
$ws.Range("H86").Value = 5806.0586
$ws.Range("I86").Value = 5766.6
$ws.Range("J86").Value = 5862.4287
$ws.Range("K86").Value = 5766.6
$ws.Range("L86").Value = 5862.4287
$ws.Range("M86").Value = -4643.6
$ws.Range("N86").Value = -8108.4287
$ws.Range("H89").Value = 5806.0586
$ws.Range("I89").Value = 5766.6
$ws.Range("J89").Value = 5862.4287
$ws.Range("K89").Value = 28833
$ws.Range("L89").Value = 29312.1435
$ws.Range("M89").Value = -23217
$ws.Range("N89").Value = -40544.14350000001
$ws.Range("H99").Value = 2670.7693
$ws.Range("I99").Value = 2602.4
$ws.Range("J99").Value = 2713.5
$ws.Range("K99").Value = 2602.4
$ws.Range("L99").Value = 2713.5
$ws.Range("M99").Value = -1104.4
$ws.Range("N99").Value = -5709.5
$ws.Range("H101").Value = 29999.5
$ws.Range("J101").Value = 29999.5
$ws.Range("L101").Value = 29999.5
$ws.Range("N101").Value = -36489.5
$ws.Range("H126").Value = 2670.7693
$ws.Range("I126").Value = 2602.4
$ws.Range("J126").Value = 2713.5
$ws.Range("K126").Value = 7807.200000000001
$ws.Range("L126").Value = 8140.5
$ws.Range("M126").Value = -5337.200000000001
$ws.Range("N126").Value = -13080.5
$ws.Range("H132").Value = 2337.0527
$ws.Range("I132").Value = 2145.484
$ws.Range("K132").Value = 6436.451999999999
$ws.Range("M132").Value = -3906.451999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 69.818184
$ws.Range("I2").Value = 77.5
$ws.Range("K2").Value = 465
$ws.Range("M2").Value = -352
$ws.Range("H11").Value = 69215.37
$ws.Range("I11").Value = 61386.35
$ws.Range("K11").Value = 184159.05
$ws.Range("M11").Value = -184019.05
$ws.Range("H26").Value = 189.8
$ws.Range("I26").Value = 175
$ws.Range("K26").Value = 525
$ws.Range("M26").Value = -237
$ws.Range("H46").Value = 204
$ws.Range("J46").Value = 399
$ws.Range("L46").Value = 1197
$ws.Range("N46").Value = -1379
$ws.Range("H92").Value = 380.5
$ws.Range("I92").Value = 220.57143
$ws.Range("K92").Value = 661.71429
$ws.Range("M92").Value = 586.28571
$ws.Range("H134").Value = 2999
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2353.6
$ws.Range("J80").Value = 2419.4443
$ws.Range("L80").Value = 2419.4443
$ws.Range("N80").Value = -4415.4443
$ws.Range("H83").Value = 2353.6
$ws.Range("J83").Value = 2419.4443
$ws.Range("L83").Value = 12097.2215
$ws.Range("N83").Value = -22081.2215
$ws.Range("H113").Value = 3387.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 12776.667
$ws.Range("J136").Value = 12776.667
$ws.Range("L136").Value = 38330.001
$ws.Range("N136").Value = -43430.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 626.4286
$ws.Range("I107").Value = 522.7273
$ws.Range("K107").Value = 1568.1819
$ws.Range("M107").Value = 351.8181
$ws.Range("H113").Value = 7138.857
$ws.Range("J113").Value = 10361.5
$ws.Range("L113").Value = 31084.5
$ws.Range("N113").Value = -35424.5
$ws.Range("H122").Value = 3788
$ws.Range("I122").Value = 2082
$ws.Range("K122").Value = 6246
$ws.Range("M122").Value = -3796
$ws.Range("H126").Value = 1505.25
$ws.Range("I126").Value = 1449.8572
$ws.Range("K126").Value = 4349.571599999999
$ws.Range("M126").Value = -1879.571599999999
$ws.Range("H132").Value = 649772.4399999999
$ws.Range("I132").Value = 913387.5600000001
$ws.Range("K132").Value = 2740162.68
$ws.Range("M132").Value = -2737632.68
